# The commit updates the "Login" test-data sheet so that the password
# cell (B2) no longer shows the old plaintext password ("#123ew4567dsf8w")
# but instead shows that the password has been changed.
#
# Target cell: B2 ("Password" row for sbuda@gmail.com) goes from
#   "#123ew4567dsf8w"  ->  "PasswordChanged"
#
# The leading apostrophe forces Excel's "stored as text / quote-prefixed"
# entry semantics, which keeps the cell's existing number format / style
# (border only, quote-prefixed text) instead of resetting it when the new
# value is written - matching the original cell's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

$ws.Cells.Item(2, 2).Value = "'PasswordChanged"

$wb.Save()
